$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-10-23 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-10-24 Friday", 2) | Out-Null
$d.Content.Find.Execute("85÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "75÷3=", 2) | Out-Null
$d.Content.Find.Execute("24÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "25÷5=", 2) | Out-Null
$d.Content.Find.Execute("38÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "43÷3=", 2) | Out-Null
$d.Content.Find.Execute("61÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "10÷2=", 2) | Out-Null
$d.Content.Find.Execute("94÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "33÷6=", 2) | Out-Null
$d.Content.Find.Execute("12÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "49÷5=", 2) | Out-Null
$d.Content.Find.Execute("98÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "52÷7=", 2) | Out-Null
$d.Content.Find.Execute("47÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "11÷8=", 2) | Out-Null
$d.Content.Find.Execute("32÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "80÷6=", 2) | Out-Null
$d.Content.Find.Execute("81÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "20÷6=", 2) | Out-Null
$d.Content.Find.Execute("36÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "71÷6=", 2) | Out-Null
$d.Content.Find.Execute("42÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "30÷8=", 2) | Out-Null
$d.Content.Find.Execute("83÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "35÷3=", 2) | Out-Null
$d.Content.Find.Execute("66÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "15÷4=", 2) | Out-Null
$d.Content.Find.Execute("10÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "21÷8=", 2) | Out-Null
$d.Content.Find.Execute("50÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "78÷2=", 2) | Out-Null
$d.Content.Find.Execute("90÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "69÷5=", 2) | Out-Null
$d.Content.Find.Execute("67÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "87÷6=", 2) | Out-Null
$d.Content.Find.Execute("26÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "71÷5=", 2) | Out-Null
$d.Content.Find.Execute("24÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "35÷9=", 2) | Out-Null
$d.Content.Find.Execute("68÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "62÷2=", 2) | Out-Null
$d.Content.Find.Execute("42÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "89÷2=", 2) | Out-Null
$d.Content.Find.Execute("19÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "72÷5=", 2) | Out-Null
$d.Content.Find.Execute("92÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "67÷9=", 2) | Out-Null
$d.Content.Find.Execute("91÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "25÷9=", 2) | Out-Null
